$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 281
$ws.Range("F6").Value = 663
$ws.Range("F7").Value = 21
$ws.Range("F8").Value = 2704
$ws.Range("F10").Value = 6363
$ws.Range("F11").Value = 2401
$ws.Range("F13").Value = 24
$ws.Range("F15").Value = 2569
$ws.Range("F16").Value = 30
$ws.Range("F17").Value = 25
$ws.Range("F18").Value = 6793
$ws.Range("F19").Value = 249
$ws.Range("F20").Value = 88
$ws.Range("F21").Value = 184
$ws.Range("F22").Value = 109
$ws.Range("F24").Value = 7666
$ws.Range("F32").Value = 50
$ws.Range("F39").Value = 38
$ws.Range("F44").Value = 595
$ws.Range("F45").Value = 3601
$ws.Range("F46").Value = 130
$ws.Range("F47").Value = 1156
$ws.Range("F48").Value = 92

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 233
$ws.Range("F7").Value = 104
$ws.Range("F8").Value = 24
$ws.Range("F16").Value = 18

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 281
$ws.Range("F4").Value = 663
$ws.Range("F5").Value = 21
$ws.Range("F7").Value = 2704
$ws.Range("F9").Value = 233
$ws.Range("F10").Value = 6363
$ws.Range("F11").Value = 2401
$ws.Range("F12").Value = 104
$ws.Range("F13").Value = 24
$ws.Range("F15").Value = 2569
$ws.Range("F16").Value = 30
$ws.Range("F17").Value = 24
$ws.Range("F19").Value = 25
$ws.Range("F20").Value = 6793
$ws.Range("F21").Value = 249
$ws.Range("F22").Value = 88
$ws.Range("F23").Value = 184
$ws.Range("F25").Value = 7667
$ws.Range("F38").Value = 2572
$ws.Range("F39").Value = 38
$ws.Range("F43").Value = 595
$ws.Range("F45").Value = 3601
$ws.Range("F46").Value = 130
$ws.Range("F47").Value = 18
$ws.Range("F48").Value = 1156
$ws.Range("F49").Value = 92
